# Update the Metadata sheet of the ConceptMap workbook:
#  - URL now points at the new /ig/ path
#  - Date bumped to the new publication timestamp
#  - Jurisdiction value filled in with "FRANCE" (was blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ConceptMap/PN13-FHIR-prescmed-patient-sans-INS-conceptmap"
$ws.Range("B8").Value = "2026-01-15T08:54:26+00:00"
$ws.Range("B11").Value = "FRANCE"
